$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: client 8888 - "DESPENSA X" / "CENTRO MBARE" / "NOMBRE VENDEDOR EJEMPLO" / "JU"
$ws.Range("A2").Value = 8888
$ws.Range("B2").Value = "DESPENSA X"
$ws.Range("C2").Value = "CENTRO MBARE"
$ws.Range("F2").Value = "NOMBRE VENDEDOR EJEMPLO"
$ws.Range("H2").Value = "JU"

# Row 3: client 9999 - "BODEGA EJEMPLO" / "CENTRO MBARE" / "NOMBRE VENDEDOR EJEMPLO" / "JU"
$ws.Range("A3").Value = 9999
$ws.Range("B3").Value = "BODEGA EJEMPLO"
$ws.Range("C3").Value = "CENTRO MBARE"
$ws.Range("F3").Value = "NOMBRE VENDEDOR EJEMPLO"
$ws.Range("H3").Value = "JU"

# Update the active selection to A3 (as seen in the saved file)
$ws.Range("A3").Select()
